$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric results for rows 2-4 (random_forest, lsboost, neural_network)
# Columns: B=RMSE, C=NRMSE, D=MAE, E=RSE, F=RRSE, G=RAE, H=R2, I=Corr Coeff

$ws.Range("B2").Value = 3.6828651750473793
$ws.Range("C2").Value = 0.2309583152706009
$ws.Range("D2").Value = 2.8038281809422436
$ws.Range("E2").Value = 0.38833261787398821
$ws.Range("F2").Value = 0.6231633958072218
$ws.Range("G2").Value = 0.58644644003489277
$ws.Range("H2").Value = 0.61166738212601179
$ws.Range("I2").Value = 0.78429223198529008

$ws.Range("B3").Value = 3.8552902164632212
$ws.Range("C3").Value = 0.24177136304266714
$ws.Range("D3").Value = 2.9545446798829547
$ws.Range("E3").Value = 0.42554587252688642
$ws.Range("F3").Value = 0.65233877128903384
$ws.Range("G3").Value = 0.61797018134653026
$ws.Range("H3").Value = 0.57445412747311364
$ws.Range("I3").Value = 0.76179211881353304

$ws.Range("B4").Value = 3.9426863240502232
$ws.Range("C4").Value = 0.24725211153877316
$ws.Range("D4").Value = 2.9479807549391852
$ws.Range("E4").Value = 0.44505807280725224
$ws.Range("F4").Value = 0.66712672919562455
$ws.Range("G4").Value = 0.61659727610144632
$ws.Range("H4").Value = 0.55494192719274782
$ws.Range("I4").Value = 0.74520754331358474
